$d = $word.ActiveDocument

$replacements = @(
    @("58×43=", "36×23="),
    @("98×59=", "25×93="),
    @("85×27=", "52×75="),
    @("37×20=", "83×41="),
    @("69×11=", "31×36="),
    @("25×39=", "80×20="),
    @("84×26=", "46×81="),
    @("58×59=", "18×28="),
    @("24×99=", "18×67="),
    @("86×17=", "72×73="),
    @("73×90=", "71×97="),
    @("51×15=", "92×37="),
    @("79×91=", "29×19="),
    @("29×26=", "60×76="),
    @("34×50=", "48×53="),
    @("24×12=", "75×50="),
    @("62×91=", "73×37="),
    @("68×59=", "23×30="),
    @("78×79=", "42×64="),
    @("88×41=", "65×65="),
    @("98×42=", "99×89="),
    @("84×39=", "84×54="),
    @("41×91=", "88×37="),
    @("63×64=", "65×34="),
    @("96×48=", "25×14=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
